# Insert a new weekly data row at row 140 (Vega Central Mapocho de Santiago - Haba),
# shifting existing rows 140-166 down to 141-167, then populate the new row 140
# with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 140; this shifts rows 140:166 down to 141:167
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new record's data
$ws.Cells.Item(140, 1).Value = 9
$ws.Cells.Item(140, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(140, 3).Value = "Metropolitana"
$ws.Cells.Item(140, 4).Value = 44504
$ws.Cells.Item(140, 5).Value = 13
$ws.Cells.Item(140, 6).Value = 100112026
$ws.Cells.Item(140, 7).Value = "Haba"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 79
$ws.Cells.Item(140, 11).Value = 7000
$ws.Cells.Item(140, 12).Value = 8000
$ws.Cells.Item(140, 13).Value = 7506
$ws.Cells.Item(140, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(140, 15).Value = "Región Metropolitana"
$ws.Cells.Item(140, 16).Value = 300
$ws.Cells.Item(140, 17).Value = 25
$ws.Cells.Item(140, 18).Value = "Hortaliza"

# Apply the same date format style used by the other rows in column D
$ws.Cells.Item(140, 4).NumberFormat = $ws.Cells.Item(141, 4).NumberFormat
